# Updates cryptos list prices (D) and volume-change percentages (E)
# Values are plain text (not numbers), so we force Text format before
# assigning, then restore the default "Normal" style so the cell keeps
# the same (unstyled) appearance as before - only its text changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.398.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.824.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07881"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9632"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.831.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.869"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.079"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06596"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001022"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.391.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.304"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.056.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.273"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09310"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9354"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.575"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.224"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.321"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.113"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "

$ws.Range("E40").Value = "  -0.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5760"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1821"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.958"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5418"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.872"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06564"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "109.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.86%  "

$ws.Range("E51").Value = "  -33.20%  "
